$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E1").Value = "DMV_trein_attempts"
$ws.Range("E2:E127").Formula = "=D2/90.9*100"
$ws.Range("E2:E127").NumberFormat = "0"
$null = $ws.Range("E1").Select()
$ws.PageSetup.Orientation = 1
